# Applies the "Updated cryptos list" data refresh (Wed Apr 26 16:57:44 UTC 2023).
# Column D = Price, Column E = Volume(1h) change; rows 47/48 (EnergySwap <-> Decentraland) swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.802.15"
$ws.Range("E2").Value = "  +8.42%  "
$ws.Range("D3").Value = "1.953.28"
$ws.Range("E3").Value = "  +6.77%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "342.79"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.4779"
$ws.Range("E7").Value = "  +4.46%  "
$ws.Range("D8").Value = "0.4145"
$ws.Range("E8").Value = "  +8.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.00"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("D10").Value = "0.08259"
$ws.Range("E10").Value = "  +5.00%  "
$ws.Range("D11").Value = "1.037"
$ws.Range("E11").Value = "  +8.16%  "
$ws.Range("D12").Value = "22.74"
$ws.Range("E12").Value = "  +7.91%  "
$ws.Range("D13").Value = "1.939.66"
$ws.Range("E13").Value = "  +5.93%  "
$ws.Range("D14").Value = "6.189"
$ws.Range("E14").Value = "  +5.94%  "
$ws.Range("D15").Value = "7.422"
$ws.Range("E15").Value = "  +5.21%  "
$ws.Range("D16").Value = "92.25"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("D19").Value = "0.06706"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "18.06"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "29.758.40"
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("D23").Value = "5.585"
$ws.Range("E23").Value = "  +5.39%  "
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "2.180.08"
$ws.Range("E26").Value = "  +6.58%  "
$ws.Range("D27").Value = "161.88"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "20.22"
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").Value = "2.186"
$ws.Range("E29").Value = "  +6.87%  "
$ws.Range("D30").Value = "5.706"
$ws.Range("E30").Value = "  +7.87%  "
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("E32").Value = "  +8.71%  "
$ws.Range("D33").Value = "0.09634"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "1.481"
$ws.Range("E34").Value = "  +12.47%  "
$ws.Range("D35").Value = "3.688"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").Value = "5.525"
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("D37").Value = "0.06298"
$ws.Range("E37").Value = "  +6.51%  "
$ws.Range("D38").Value = "0.02322"
$ws.Range("D39").Value = "8.502"
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.190"
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("D41").Value = "0.6106"
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("D42").Value = "10.75"
$ws.Range("E42").Value = "  +8.51%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "0.1896"
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("D45").Value = "2.381"
$ws.Range("E45").Value = "  +32.60%  "
$ws.Range("D46").Value = "1.261"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5724"
$ws.Range("E47").Value = "  +6.04%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "12.48"
$ws.Range("E48").Value = "  +6.30%  "
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").Value = "0.07354"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("D51").Value = "113.73"
$ws.Range("E51").Value = "  +2.95%  "
